$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O11").Value = 68001.31
$ws.Range("O12").Value = 868.02
$ws.Range("O15").Value = 2058.54
$ws.Range("N23").Value = 192424.74
$ws.Range("O23").Value = 191833.83
